$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new county name "Bala" to cell A2 (as a shared string)
$ws.Range("A2").Value = "Bala"

# Move the active selection to A3, matching the post-edit state
$ws.Range("A3").Select()
